$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 53.08554433333333
$ws.Range("H2").Value = 159.256633
$ws.Range("I2").Value = 0.5980804157037442
$ws.Range("J2").Value = 0.598080415703744
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 6213.374338112126
$ws.Range("R2").Value = 55920.36904300914
$ws.Range("S2").Value = 0.1940989634328117
$ws.Range("T2").Value = 0.1940989634328117

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 53.08554433333333
$ws.Range("H3").Value = 159.256633
$ws.Range("I3").Value = 0.5980804157037442
$ws.Range("J3").Value = 0.598080415703744
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 5392.431575240322
$ws.Range("R3").Value = 48531.88417716289
$ws.Range("S3").Value = 0.1684536166952609
$ws.Range("T3").Value = 0.1684536166952608

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 53.08554433333333
$ws.Range("H4").Value = 159.256633
$ws.Range("I4").Value = 0.5980804157037442
$ws.Range("J4").Value = 0.598080415703744
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 7539.56942167566
$ws.Range("R4").Value = 67856.12479508095
$ws.Range("S4").Value = 0.2355278355756715
$ws.Range("T4").Value = 0.2355278355756714

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 30.06295833333333
$ws.Range("H5").Value = 90.188875
$ws.Range("I5").Value = 0.3386998634578254
$ws.Range("J5").Value = 0.3386998634578254
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 3518.705820612208
$ws.Range("R5").Value = 31668.35238550988
$ws.Range("S5").Value = 0.1099204900977118
$ws.Range("T5").Value = 0.1099204900977118

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 30.06295833333333
$ws.Range("H6").Value = 90.188875
$ws.Range("I6").Value = 0.3386998634578254
$ws.Range("J6").Value = 0.3386998634578254
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 3053.796429850444
$ws.Range("R6").Value = 27484.167868654
$ws.Range("S6").Value = 0.09539723334114944
$ws.Range("T6").Value = 0.09539723334114943

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 30.06295833333333
$ws.Range("H7").Value = 90.188875
$ws.Range("I7").Value = 0.3386998634578254
$ws.Range("J7").Value = 0.3386998634578254
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 4269.745449945111
$ws.Range("R7").Value = 38427.709049506
$ws.Range("S7").Value = 0.1333821400189642
$ws.Range("T7").Value = 0.1333821400189642

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.611374666666666
$ws.Range("H8").Value = 16.834124
$ws.Range("I8").Value = 0.06321972083843048
$ws.Range("J8").Value = 0.06321972083843048
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 656.7808956892707
$ws.Range("R8").Value = 5911.028061203436
$ws.Range("S8").Value = 0.02051711101225792
$ws.Range("T8").Value = 0.02051711101225792

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.611374666666666
$ws.Range("H9").Value = 16.834124
$ws.Range("I9").Value = 0.06321972083843048
$ws.Range("J9").Value = 0.06321972083843048
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 570.0036481313209
$ws.Range("R9").Value = 5130.032833181888
$ws.Range("S9").Value = 0.01780628547946567
$ws.Range("T9").Value = 0.01780628547946566

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.611374666666666
$ws.Range("H10").Value = 16.834124
$ws.Range("I10").Value = 0.06321972083843048
$ws.Range("J10").Value = 0.06321972083843048
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 796.9655276530702
$ws.Range("R10").Value = 7172.689748877632
$ws.Range("S10").Value = 0.0248963243467069
$ws.Range("T10").Value = 0.0248963243467069
